$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 3481.6365
$ws.Range("I101").Value = 432.4
$ws.Range("K101").Value = 1297.2
$ws.Range("M101").Value = 324.8000000000002
$ws.Range("H112").Value = 3809.8518
$ws.Range("J112").Value = 3809.8518
$ws.Range("L112").Value = 11429.5554
$ws.Range("N112").Value = -13645.5554
$ws.Range("H132").Value = 4210.143
$ws.Range("I132").Value = 1395.7
$ws.Range("K132").Value = 4187.1
$ws.Range("M132").Value = -1657.1
$ws.Range("H137").Value = 26413354
$ws.Range("I137").Value = 1003524.9
$ws.Range("K137").Value = 3010574.7
$ws.Range("M137").Value = -3008024.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26228.121
$ws.Range("I32").Value = 26225.275
$ws.Range("K32").Value = 26225.275
$ws.Range("M32").Value = -25938.275
$ws.Range("H61").Value = 4340.364
$ws.Range("I61").Value = 4340.364
$ws.Range("K61").Value = 4340.364
$ws.Range("M61").Value = -4128.364
$ws.Range("H74").Value = 27780572
$ws.Range("I74").Value = 31252830
$ws.Range("K74").Value = 31252830
$ws.Range("M74").Value = -31251956
$ws.Range("H77").Value = 27780572
$ws.Range("I77").Value = 31252830
$ws.Range("K77").Value = 156264150
$ws.Range("M77").Value = -156259782
$ws.Range("H104").Value = 44931.25
$ws.Range("J104").Value = 44931.25
$ws.Range("L104").Value = 44931.25
$ws.Range("N104").Value = -51919.25
$ws.Range("H122").Value = 3884.9534
$ws.Range("J122").Value = 7346.5
$ws.Range("L122").Value = 22039.5
$ws.Range("N122").Value = -26939.5
$ws.Range("H132").Value = 2677
$ws.Range("I132").Value = 2677
$ws.Range("K132").Value = 8031
$ws.Range("M132").Value = -5501
$ws.Range("H133").Value = 81666.664
$ws.Range("J133").Value = 81666.664
$ws.Range("L133").Value = 81666.664
$ws.Range("N133").Value = -86726.664
$ws.Range("H136").Value = 4340.364
$ws.Range("I136").Value = 4340.364
$ws.Range("K136").Value = 13021.092
$ws.Range("M136").Value = -10471.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2536.2122
$ws.Range("I107").Value = 2415.6956
$ws.Range("K107").Value = 2415.6956
$ws.Range("M107").Value = -495.6956
$ws.Range("H134").Value = 3227.76
$ws.Range("I134").Value = 2563.2727
$ws.Range("K134").Value = 7689.8181
$ws.Range("M134").Value = -5154.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 68580.5
$ws.Range("J18").Value = 68580.5
$ws.Range("L18").Value = 68580.5
$ws.Range("N18").Value = -69040.5
$ws.Range("H31").Value = 20002520
$ws.Range("I31").Value = 21278876
$ws.Range("K31").Value = 21278876
$ws.Range("M31").Value = -21278581
$ws.Range("H34").Value = 20002520
$ws.Range("I34").Value = 21278876
$ws.Range("K34").Value = 21278876
$ws.Range("M34").Value = -21278674
$ws.Range("H58").Value = 2644.5
$ws.Range("I58").Value = 2644.5
$ws.Range("K58").Value = 2644.5
$ws.Range("M58").Value = -2441.5
$ws.Range("H136").Value = 2644.5
$ws.Range("I136").Value = 2644.5
$ws.Range("K136").Value = 7933.5
$ws.Range("M136").Value = -5383.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5089283
$ws.Range("I4").Value = 136157.38
$ws.Range("J4").Value = 11693450
$ws.Range("K4").Value = 408472.14
$ws.Range("L4").Value = 35080350
$ws.Range("M4").Value = -408360.14
$ws.Range("N4").Value = -35080574
$ws.Range("H6").Value = 84668.164
$ws.Range("I6").Value = 111557.555
$ws.Range("K6").Value = 334672.665
$ws.Range("M6").Value = -334559.665
$ws.Range("H107").Value = 390.8846
$ws.Range("J107").Value = 454.72223
$ws.Range("L107").Value = 1364.16669
$ws.Range("N107").Value = -5204.16669
$ws.Range("H132").Value = 1173.1714
$ws.Range("I132").Value = 981.5185
$ws.Range("J132").Value = 1820
$ws.Range("K132").Value = 8833.666499999999
$ws.Range("L132").Value = 16380
$ws.Range("M132").Value = -6303.666499999999
$ws.Range("N132").Value = -21440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7702.68
$ws.Range("I70").Value = 7857.0586
$ws.Range("K70").Value = 7857.0586
$ws.Range("M70").Value = -7587.0586
$ws.Range("H73").Value = 7702.68
$ws.Range("I73").Value = 7857.0586
$ws.Range("K73").Value = 7857.0586
$ws.Range("M73").Value = -6921.0586
$ws.Range("H102").Value = 10873242
$ws.Range("J102").Value = 4758
$ws.Range("L102").Value = 4758
$ws.Range("N102").Value = -8002
$ws.Range("H113").Value = 775
$ws.Range("I113").Value = 775
$ws.Range("K113").Value = 775
$ws.Range("M113").Value = 1395
$ws.Range("H122").Value = 404827.62
$ws.Range("I122").Value = 835391.9399999999
$ws.Range("K122").Value = 2506175.82
$ws.Range("M122").Value = -2503725.82
$ws.Range("H132").Value = 289746.72
$ws.Range("I132").Value = 401685.4
$ws.Range("J132").Value = 9900
$ws.Range("K132").Value = 1205056.2
$ws.Range("L132").Value = 29700
$ws.Range("M132").Value = -1202526.2
$ws.Range("N132").Value = -34760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 25000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H46").Value = 3071.1428
$ws.Range("J46").Value = 3082.8333
$ws.Range("L46").Value = 3082.8333
$ws.Range("N46").Value = -3458.8333
$ws.Range("H61").Value = 7950
$ws.Range("I61").Value = 7950
$ws.Range("K61").Value = 7950
$ws.Range("M61").Value = -7748
$ws.Range("H100").Value = 2711
$ws.Range("I100").Value = 2496
$ws.Range("K100").Value = 2496
$ws.Range("M100").Value = -1955
$ws.Range("H113").Value = 7950
$ws.Range("I113").Value = 7950
$ws.Range("K113").Value = 7950
$ws.Range("M113").Value = -5780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9333.143
$ws.Range("J74").Value = 9333.143
$ws.Range("L74").Value = 9333.143
$ws.Range("N74").Value = -11205.143
$ws.Range("H77").Value = 9333.143
$ws.Range("J77").Value = 9333.143
$ws.Range("L77").Value = 27999.429
$ws.Range("N77").Value = -37359.429
$ws.Range("H124").Value = 10000
$ws.Range("J124").Value = 10000
$ws.Range("L124").Value = 10000
$ws.Range("N124").Value = -19820
$ws.Range("H132").Value = 1492.1216
$ws.Range("I132").Value = 1456.6666
$ws.Range("J132").Value = 1498.9839
$ws.Range("K132").Value = 4369.9998
$ws.Range("L132").Value = 4496.9517
$ws.Range("M132").Value = -1839.9998
$ws.Range("N132").Value = -9556.9517
$ws.Range("H136").Value = 3541.3057
$ws.Range("I136").Value = 2451.3572
$ws.Range("K136").Value = 7354.071599999999
$ws.Range("M136").Value = -4804.071599999999
